$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.869.35"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.771.09"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'355.24"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'108.50"
$ws.Range("E6").Value = "  -4.84%  "
$ws.Range("D7").Value = "'0.560"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.587"
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").Value = "'39.94"
$ws.Range("E10").Value = "  -5.09%  "
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "'19.30"
$ws.Range("E13").Value = "  -4.19%  "
$ws.Range("D14").Value = "'7.56"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "3.208.14"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "2.775.43"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "'0.929"
$ws.Range("E17").Value = "  +3.34%  "
$ws.Range("D18").Value = "51.742.70"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'12.98"
$ws.Range("E21").Value = "  -4.86%  "
$ws.Range("D22").Value = "0.0₃0971"
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").Value = "'273.87"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "'69.57"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "'2.72"
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("D26").Value = "'26.44"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "'10.09"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("B31").Value = "OKB"
$ws.Range("C31").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D31").Value = "'51.36"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("B32").Value = "VeChain"
$ws.Range("C32").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D32").Value = "'0.0460"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "'33.58"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").Value = "'5.68"
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("D35").Value = "'5.29"
$ws.Range("E35").Value = "  +9.34%  "
$ws.Range("D36").Value = "'0.0832"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "'3.18"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("D39").Value = "'18.07"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "'1.98"
$ws.Range("E40").Value = "  -5.29%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.52"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.114"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "'121.85"
$ws.Range("E43").Value = "  -4.99%  "
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").Value = "'21.65"
$ws.Range("E45").Value = "  -8.51%  "
$ws.Range("D46").Value = "2.051.92"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "'3.23"
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "'5.67"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "'0.919"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").Value = "'8.84"
$ws.Range("E51").Value = "  -0.64%  "
